$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("C2").Value = 39
$ws.Range("F2").Value = 51
$ws.Range("I3").Value = 111
$ws.Range("B6").Value = 225
$ws.Range("C6").Value = 280
$ws.Range("D6").Value = 251
$ws.Range("E6").Value = 250
$ws.Range("F6").Value = 326
$ws.Range("G6").Value = 290
$ws.Range("H6").Value = 252
$ws.Range("B7").Value = 300
$ws.Range("C7").Value = 376
$ws.Range("D7").Value = 392
$ws.Range("E7").Value = 382
$ws.Range("F7").Value = 464
$ws.Range("G7").Value = 420
$ws.Range("H7").Value = 386
$ws.Range("I7").Value = 505

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("D2").Value = 2
$ws.Range("H5").Value = 3
$ws.Range("B8").Value = 21
$ws.Range("F8").Value = 33
$ws.Range("E18").Value = 11
$ws.Range("H26").Value = 10
$ws.Range("F27").Value = 34
$ws.Range("F28").Value = 9
$ws.Range("E31").Value = 31
$ws.Range("C35").Value = 21
$ws.Range("F35").Value = 29
$ws.Range("C46").Value = 19
$ws.Range("F52").Value = 46
$ws.Range("I53").Value = 5
$ws.Range("F71").Value = 7
$ws.Range("G73").Value = 9
$ws.Range("I73").Value = 11
$ws.Range("G75").Value = 9
$ws.Range("C76").Value = 14
$ws.Range("E76").Value = 15
$ws.Range("H76").Value = 14
$ws.Range("B79").Value = 10
$ws.Range("C91").Value = 3
$ws.Range("E94").Value = 4
$ws.Range("B97").Value = 300
$ws.Range("C97").Value = 376
$ws.Range("D97").Value = 392
$ws.Range("E97").Value = 382
$ws.Range("F97").Value = 464
$ws.Range("G97").Value = 420
$ws.Range("H97").Value = 386
$ws.Range("I97").Value = 505

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 9

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("C6").Value = 10
$ws.Range("E6").Value = 9
$ws.Range("H6").Value = 10
$ws.Range("C7").Value = 14
$ws.Range("E7").Value = 15
$ws.Range("H7").Value = 14

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("B5").Value = 13
$ws.Range("F5").Value = 22
$ws.Range("B6").Value = 21
$ws.Range("F6").Value = 33

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("E6").Value = 22
$ws.Range("E7").Value = 31

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("C6").Value = 18
$ws.Range("F6").Value = 11
$ws.Range("C7").Value = 21
$ws.Range("F7").Value = 29

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("H5").Value = 3
$ws.Range("H6").Value = 3

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("B4").Value = 9
$ws.Range("B5").Value = 10

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("F6").Value = 25
$ws.Range("F7").Value = 34

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("F2").Value = 4
$ws.Range("F7").Value = 46

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 3

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I2").Value = 2
$ws.Range("G4").Value = 8
$ws.Range("G5").Value = 9
$ws.Range("I5").Value = 11

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("F5").Value = 9
$ws.Range("F6").Value = 9

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("F4").Value = 6
$ws.Range("F5").Value = 7

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("H4").Value = 10
$ws.Range("H5").Value = 10

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("E5").Value = 10
$ws.Range("E6").Value = 11

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 2

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I3").Value = 3
$ws.Range("I5").Value = 5

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("C2").Value = 4
$ws.Range("C6").Value = 19

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 4
